$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before the old
#    "2022-Q2" sheet, which is currently the 2nd sheet).
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $zj)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the "2022-Q3" sheet with the fund holding table.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3Data = @(
  ,@(0, "506000", "南方科创板 3 年定期开放混合", "25.21", "79.80", "3.65", "0.9202", 8)
  ,@(1, "001150", "融通互联网传媒灵活配置混合", "8.61", "88.88", "5.41", "0.4658", 1)
  ,@(2, "006348", "银华盛利混合A", "12.76", "86.66", "3.53", "0.4504", 5)
  ,@(3, "001227", "中邮信息产业灵活配置混合", "5.91", "85.43", "4.01", "0.2370", 2)
  ,@(4, "009715", "汇添富策略增长灵活配置混合", "3.87", "88.49", "4.89", "0.1892", 6)
  ,@(5, "001541", "汇添富民营新动力股票", "2.21", "91.50", "4.87", "0.1076", 3)
  ,@(6, "013641", "博道成长智航股票A", "10.27", "90.29", "0.98", "0.1006", 3)
  ,@(7, "506008", "长城科创两年定开混合A", "3.12", "78.09", "2.95", "0.0920", 6)
  ,@(8, "398011", "中海分红增利混合", "2.50", "91.65", "3.43", "0.0858", 7)
  ,@(9, "013369", "汇添富自主核心科技一年持有混合A", "2.47", "65.13", "3.08", "0.0761", 6)
  ,@(10, "013642", "博道成长智航股票C", "7.24", "90.29", "0.98", "0.0710", 3)
  ,@(11, "001728", "银华战略新兴灵活配置定期开放混合", "1.39", "97.07", "4.35", "0.0605", 3)
  ,@(12, "008602", "方正富邦新兴成长混合A", "1.23", "86.03", "4.38", "0.0539", 3)
  ,@(13, "001275", "中邮创新优势灵活配置混合", "1.04", "86.62", "3.94", "0.0410", 3)
  ,@(14, "015684", "银华盛利混合C", "1.05", "86.66", "3.53", "0.0371", 5)
  ,@(15, "002213", "中海顺鑫灵活配置混合", "0.75", "91.58", "3.46", "0.0260", 6)
  ,@(16, "013370", "汇添富自主核心科技一年持有混合C", "0.67", "65.13", "3.08", "0.0206", 6)
  ,@(17, "003659", "山西证券策略精选灵活配置混合", "0.27", "78.35", "3.26", "0.0088", 2)
  ,@(18, "012793", "长城科创两年定开混合C", "0.11", "78.09", "2.95", "0.0032", 6)
  ,@(19, "008603", "方正富邦新兴成长混合C", "0.03", "86.03", "4.38", "0.0013", 3)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]

    # Column A: numeric running index.
    $q3.Cells.Item($r, 1).Value = $row[0]

    # Column B: fund code (kept as text even though it is all digits).
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 2).ClearFormats()

    # Column C: fund name (plain text).
    $q3.Cells.Item($r, 3).Value = $row[2]

    # Columns D-G: numeric-looking values stored as text (quote-prefix trick,
    # then drop the resulting quote-prefix style so no extra formatting is
    # introduced).
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 4).ClearFormats()
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 5).ClearFormats()
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 6).ClearFormats()
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 7).ClearFormats()

    # Column H: numeric rank.
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 3. Match formatting: header row + column A use the bold/border/centered
#    style already present on the other quarter sheets (style index 2 in the
#    original file). Copy it over from the analogous cells on the "总计"
#    sheet / the old "2022-Q2" sheet (now sheet index 3).
# ---------------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(3)

$oldQ2.Range("B1:H1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$oldQ2.Range("A2").Copy() | Out-Null
$q3.Range("A2:A21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above
#    the existing quarters, shifting 2022-Q2 / 2022-Q1 / 2021-Q4 down by one
#    row, and renumber the index column.
# ---------------------------------------------------------------------------
$zj.Range("A4").Copy() | Out-Null
$zj.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$zj.Range("A3:A5").Value = $zj.Range("A2:A4").Value
$zj.Range("B3:D5").Value = $zj.Range("B2:D4").Value

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q3"
$zj.Cells.Item(2, 3).Value = 20
$zj.Cells.Item(2, 4).Value = 3.05

$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2022-Q2"
$zj.Cells.Item(3, 3).Value = 35
$zj.Cells.Item(3, 4).Value = 5.3

$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2022-Q1"
$zj.Cells.Item(4, 3).Value = 26
$zj.Cells.Item(4, 4).Value = 5.45

$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(5, 2).Value = "2021-Q4"
$zj.Cells.Item(5, 3).Value = 1
$zj.Cells.Item(5, 4).Value = 0.4

Write-Output "done"
